$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the 40 new URL/status rows (rows 42-81) from the latest status-code crawl.
# Column A holds the request URL (text), column B the HTTP status code (number).
$ws.Range("A42").Value = 'https://www.allerganpro.com/co/es.html'
$ws.Range("B42").Value = 200
$ws.Range("A43").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/meta-navigation/clientlibs.min.css'
$ws.Range("B43").Value = 200
$ws.Range("A44").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/Mobile%20nav.png'
$ws.Range("B44").Value = 200
$ws.Range("A45").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro/publish-header.min.js'
$ws.Range("B45").Value = 200
$ws.Range("A46").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/allergan-pro/publish-header.min.css'
$ws.Range("B46").Value = 200
$ws.Range("A47").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/images/icons/mobile-menu-close.png'
$ws.Range("B47").Value = 200
$ws.Range("A48").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.css'
$ws.Range("B48").Value = 200
$ws.Range("A49").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/publish-footer.min.js'
$ws.Range("B49").Value = 200
$ws.Range("A50").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/image-extension/clientlibs.min.js'
$ws.Range("B50").Value = 200
$ws.Range("A51").Value = 'https://cag.abbvie.com:9999/jstag/managed/ruxitagent_A2Vfqru_10249220905100923.js'
$ws.Range("B51").Value = 200
$ws.Range("A52").Value = 'https://www.allerganpro.com/content/dam/abbvie-pro/co/abbvieprologo/AbbviePRO.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$ws.Range("B52").Value = 200
$ws.Range("A53").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_regular.woff2'
$ws.Range("B53").Value = 200
$ws.Range("A54").Value = 'https://www.allerganpro.com/etc.clientlibs/clientlibs/granite/jquery/granite/csrf.min.js'
$ws.Range("B54").Value = 200
$ws.Range("A55").Value = 'https://www.allerganpro.com/libs/granite/csrf/token.json'
$ws.Range("B55").Value = 200
$ws.Range("A56").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/meta-navigation/clientlibs.min.js'
$ws.Range("B56").Value = 200
$ws.Range("A57").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/button-link/clientlibs.min.js'
$ws.Range("B57").Value = 200
$ws.Range("A58").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/components/content/headline-text/clientlibs.min.js'
$ws.Range("B58").Value = 200
$ws.Range("A59").Value = 'https://www.allerganpro.com/content/dam/allergan-pro/colombia/home/Home_Articulos01.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$ws.Range("B59").Value = 200
$ws.Range("A60").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/hcpicon/hcpicon.ttf?q0neb3'
$ws.Range("B60").Value = 200
$ws.Range("A61").Value = 'https://www.allerganpro.com/content/allergan-pro/co/es/jcr:content/header/header_area/image-extension/item_1.coreimg.png/1682019219951-AbbviePRO.png'
$ws.Range("B61").Value = 302
$ws.Range("A62").Value = 'https://www.allerganpro.com/bin/public/abbvie-commons/hreflangs?resourcePath=/content/allergan-pro/co/es/jcr:content'
$ws.Range("B62").Value = 200
$ws.Range("A63").Value = 'https://consent.trustarc.com/v2/notice/hvz0wu'
$ws.Range("B63").Value = 200
$ws.Range("A64").Value = 'https://www.allerganpro.com/etc.clientlibs/abbvie-pro/clientlibs/assets/resources/fonts/roboto/roboto_bold.woff2'
$ws.Range("B64").Value = 200
$ws.Range("A65").Value = 'https://consent.trustarc.com/v2/asset/trustarc-logo-xs.svg'
$ws.Range("B65").Value = 200
$ws.Range("A66").Value = 'https://consent.trustarc.com/v2/asset/ic-close.svg'
$ws.Range("B66").Value = 200
$ws.Range("A67").Value = 'https://www.allerganpro.com/content/dam/allergan-pro/colombia/home/Home_Articulos02.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$ws.Range("B67").Value = 200
$ws.Range("A68").Value = 'https://consent.trustarc.com/v2/asset/latin.woff2'
$ws.Range("B68").Value = 200
$ws.Range("A69").Value = 'https://www.allerganpro.com/content/dam/allergan-pro/colombia/home/New%20banner%20Home%20Allergan%20pro.png/_jcr_content/renditions/cq5dam.web.1280.1280.png'
$ws.Range("B69").Value = 200
$ws.Range("A70").Value = 'https://consent.trustarc.com/v2/asset/16:19:49.763hvz0wu_AbbVieID-logo.png'
$ws.Range("B70").Value = 200
$ws.Range("A71").Value = 'https://consent.trustarc.com/v2/consentcategories/getnonemptyindexes?cmId=hvz0wu&referer=&fullURL=https%3A%2F%2Fwww.allerganpro.com%2Fco%2Fes.html&category='
$ws.Range("B71").Value = 200
$ws.Range("A72").Value = 'https://consent-reporting.trustarc.com/api/user-action/bannermsg?action=views&domain=hvz0wu&behavior=implied&country=bd&language=en&rand=0.25405217145246284&session=e0b1cd18-5d87-4b2c-ad3c-b046de20bccc&userType=NEW'
$ws.Range("B72").Value = 202
$ws.Range("A73").Value = 'https://consent-reporting.trustarc.com/api/user-action/log?action=impression&domain=hvz0wu&behavior=implied&country=bd&language=en&rand=0.8466606062203768&session=e0b1cd18-5d87-4b2c-ad3c-b046de20bccc&userType=NEW'
$ws.Range("B73").Value = 202
$ws.Range("A74").Value = 'https://www.allerganpro.com/content/allergan-pro/co/es/jcr%3acontent/header/header_area/image-extension/item_1.coreimg.png/1707333003504.png'
$ws.Range("B74").Value = 200
$ws.Range("A75").Value = 'https://consent.trustarc.com/v2/asset/ic-error.svg'
$ws.Range("B75").Value = 200
$ws.Range("A76").Value = 'https://consent.trustarc.com/v2/asset/ic-close-white.svg'
$ws.Range("B76").Value = 200
$ws.Range("A77").Value = 'https://www.allerganpro.com/content/allergan-pro/co/es/jcr:content/body/column_control_copy/par1-100col/column_control_18337/par1-100col/column_control/par1-5050col/image_extension_copy/item_1.coreimg.png/1637251107888-Home_Articulos01.png'
$ws.Range("B77").Value = 200
$ws.Range("A78").Value = 'https://www.allerganpro.com/content/allergan-pro/co/es/jcr:content/body/column_control_copy/par1-100col/column_control_18337/par1-100col/column_control/par2-5050col/image_extension_copy/item_1.coreimg.png/1655224770198-Home_Articulos02.png'
$ws.Range("B78").Value = 200
$ws.Range("A79").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_-2D39_sn_V1FMRLJ2C0LR8N0PRONQV2QIDHO9TGRI&svrid=-39&flavor=cors&vi=MLRPARGGUHLMQHRIUEVKEOTDAMQDALFW-0&modifiedSince=1665670355615&rf=https%3A%2F%2Fwww.allerganpro.com%2Fco%2Fes.html&bp=3&app=b90c0fbe356a6561&crc=4165458991&en=oao3vfhf&end=1'
$ws.Range("B79").Value = 200
$ws.Range("A80").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_-2D39_sn_V1FMRLJ2C0LR8N0PRONQV2QIDHO9TGRI&svrid=-39&flavor=cors&vi=MLRPARGGUHLMQHRIUEVKEOTDAMQDALFW-0&modifiedSince=1665670355615&rf=https%3A%2F%2Fwww.allerganpro.com%2Fco%2Fes.html&bp=3&app=b90c0fbe356a6561&crc=3418595279&en=oao3vfhf&end=1'
$ws.Range("B80").Value = 200
$ws.Range("A81").Value = 'https://cag.abbvie.com:9999/bf/16a183f6-c871-4082-850b-a1f7a2ecd0b1?type=js3&sn=v_4_srv_8_sn_V1FMRLJ2C0LR8N0PRONQV2QIDHO9TGRI_app-3Ab90c0fbe356a6561_1_ol_0_perc_100000_mul_1&svrid=8&flavor=cors&vi=MLRPARGGUHLMQHRIUEVKEOTDAMQDALFW-0&modifiedSince=1718839398116&rf=https%3A%2F%2Fwww.allerganpro.com%2Fco%2Fes.html&bp=3&app=b90c0fbe356a6561&crc=1379614864&en=oao3vfhf&end=1'
$ws.Range("B81").Value = 200

# Sheet view is left-to-right (matches the source report).
try { $excel.ActiveWindow.DisplayRightToLeft = $false } catch { }

# Re-assert the "numbers stored as text" ignored-error hint over the full used range.
try { $ws.Range("A1:B81").Errors.Item(9).Ignore = $true } catch { }
